$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 from "cR" to "c" (the two proximity conditions are merged into one)
$ws.Range("E1").Value = "c"

# Clear column F contents for rows 1-7 (keeps F1's existing style, just removes the value/text)
$ws.Range("F1:F7").ClearContents()

# Update row 5 values per the diff: C5=1, D5=4, E5=6 (F5 already cleared above)
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 6

# Update the active selection to match the diff (F13)
$ws.Range("F13").Select()
